# Rachel's extra data entry: append 26 more (Subject, Assignment) rows
# to sheet1, rows 41-66, continuing the existing A=subject# / B=score series.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(421, 422, 423, 424, 425, 426, 427, 428, 429, 430, 431, 432, 433, 434, 435, 436, 437, 438, 439, 440, 441, 442, 443, 444, 445, 446)
$bValues = @(9, 7, 9, 9, 4, 5, 1, 5, 3, 4, 4, 6, 7, 9, 9, 3, 1, 3, 6, 4, 9, 9, 2, 4, 7, 9)

for ($i = 0; $i -lt $aValues.Count; $i++) {
    $row = 41 + $i
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Match the new selection/scroll position recorded in the sheet view
# (active cell moves one past the last data row; view scrolled down).
$ws.Range("A70").Select()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
